$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the "Save" column values for rows 2-5
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
